$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.59333348274231
$ws.Range("B1").Value = 1.742540240287781
$ws.Range("C1").Value = 2.043889045715332
$ws.Range("D1").Value = 2.42566704750061
$ws.Range("E1").Value = 1.623193025588989
